$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamps = @{
    2 = "2025-10-17T07:09:33.201757"
    3 = "2025-10-17T07:09:33.201757"
    4 = "2025-10-17T07:09:33.201757"
    5 = "2025-10-17T07:09:33.201757"
    6 = "2025-10-17T07:09:33.202299"
    7 = "2025-10-17T07:09:33.202299"
    8 = "2025-10-17T07:09:33.202299"
    9 = "2025-10-17T07:09:33.202299"
    10 = "2025-10-17T07:09:33.202299"
    11 = "2025-10-17T07:09:33.202844"
    12 = "2025-10-17T07:09:33.202844"
    13 = "2025-10-17T07:09:33.202844"
    14 = "2025-10-17T07:09:33.202844"
    15 = "2025-10-17T07:09:33.202844"
    16 = "2025-10-17T07:09:33.203355"
    17 = "2025-10-17T07:09:33.203389"
    18 = "2025-10-17T07:09:33.203389"
    19 = "2025-10-17T07:09:33.203389"
    20 = "2025-10-17T07:09:33.203389"
    21 = "2025-10-17T07:09:33.203933"
    22 = "2025-10-17T07:09:33.203933"
    23 = "2025-10-17T07:09:33.203933"
    24 = "2025-10-17T07:09:33.203933"
    25 = "2025-10-17T07:09:33.203933"
    26 = "2025-10-17T07:09:33.203933"
    27 = "2025-10-17T07:09:33.204472"
    28 = "2025-10-17T07:09:33.204472"
    29 = "2025-10-17T07:09:33.204472"
    30 = "2025-10-17T07:09:33.204472"
    31 = "2025-10-17T07:09:33.204472"
    32 = "2025-10-17T07:09:33.205013"
    33 = "2025-10-17T07:09:33.205013"
    34 = "2025-10-17T07:09:33.205013"
    35 = "2025-10-17T07:09:33.205013"
    36 = "2025-10-17T07:09:33.205013"
    37 = "2025-10-17T07:09:33.205013"
    38 = "2025-10-17T07:09:33.205553"
    39 = "2025-10-17T07:09:33.205553"
    40 = "2025-10-17T07:09:33.205553"
    41 = "2025-10-17T07:09:33.205553"
    42 = "2025-10-17T07:09:33.205553"
    43 = "2025-10-17T07:09:33.206062"
    44 = "2025-10-17T07:09:33.206093"
    45 = "2025-10-17T07:09:33.206093"
    46 = "2025-10-17T07:09:33.266982"
    47 = "2025-10-17T07:09:33.266982"
    48 = "2025-10-17T07:09:33.267981"
    49 = "2025-10-17T07:09:33.267981"
    50 = "2025-10-17T07:09:33.267981"
    51 = "2025-10-17T07:09:33.267981"
    52 = "2025-10-17T07:09:33.267981"
    53 = "2025-10-17T07:09:33.268980"
    54 = "2025-10-17T07:09:33.268980"
    55 = "2025-10-17T07:09:33.268980"
    56 = "2025-10-17T07:09:33.268980"
    57 = "2025-10-17T07:09:33.268980"
    58 = "2025-10-17T07:09:33.268980"
    59 = "2025-10-17T07:09:33.269982"
    60 = "2025-10-17T07:09:33.269982"
    61 = "2025-10-17T07:09:33.269982"
    62 = "2025-10-17T07:09:33.269982"
    63 = "2025-10-17T07:09:33.269982"
    64 = "2025-10-17T07:09:33.270983"
    65 = "2025-10-17T07:09:33.270983"
    66 = "2025-10-17T07:09:33.270983"
    67 = "2025-10-17T07:09:33.270983"
    68 = "2025-10-17T07:09:33.270983"
    69 = "2025-10-17T07:09:33.271985"
    70 = "2025-10-17T07:09:33.271985"
    71 = "2025-10-17T07:09:33.271985"
    72 = "2025-10-17T07:09:33.271985"
    73 = "2025-10-17T07:09:33.271985"
    74 = "2025-10-17T07:09:33.272983"
    75 = "2025-10-17T07:09:33.331072"
    76 = "2025-10-17T07:09:33.331072"
    77 = "2025-10-17T07:09:33.331072"
    78 = "2025-10-17T07:09:33.331072"
    79 = "2025-10-17T07:09:33.331072"
    80 = "2025-10-17T07:09:33.331072"
    81 = "2025-10-17T07:09:33.331072"
    82 = "2025-10-17T07:09:33.331072"
    83 = "2025-10-17T07:09:33.331072"
    84 = "2025-10-17T07:09:33.331072"
    85 = "2025-10-17T07:09:33.331072"
    86 = "2025-10-17T07:09:33.331072"
    87 = "2025-10-17T07:09:33.331072"
    88 = "2025-10-17T07:09:33.331072"
    89 = "2025-10-17T07:09:33.331072"
    90 = "2025-10-17T07:09:33.331072"
    91 = "2025-10-17T07:09:33.331072"
    92 = "2025-10-17T07:09:33.331072"
    93 = "2025-10-17T07:09:33.331072"
    94 = "2025-10-17T07:09:33.331072"
    95 = "2025-10-17T07:09:33.331072"
    96 = "2025-10-17T07:09:33.331072"
    97 = "2025-10-17T07:09:33.331072"
    98 = "2025-10-17T07:09:33.331072"
    99 = "2025-10-17T07:09:33.331072"
    100 = "2025-10-17T07:09:33.331072"
    101 = "2025-10-17T07:09:33.331072"
    102 = "2025-10-17T07:09:33.331072"
    103 = "2025-10-17T07:09:33.380881"
    104 = "2025-10-17T07:09:33.380881"
    105 = "2025-10-17T07:09:33.380881"
    106 = "2025-10-17T07:09:33.380881"
    107 = "2025-10-17T07:09:33.380881"
    108 = "2025-10-17T07:09:33.380881"
    109 = "2025-10-17T07:09:33.380881"
    110 = "2025-10-17T07:09:33.380881"
    111 = "2025-10-17T07:09:33.380881"
    112 = "2025-10-17T07:09:33.380881"
}

foreach ($row in $timestamps.Keys) {
    $ws.Cells.Item($row, 26).Value = $timestamps[$row]
}
